$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update header row with short column codes
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Capitalize standalone "de", "del", "la", "los" words within names
$targets = @(
    "A22","A34","B35","B36","B46","B48","B49","B50","B51","B52",
    "B55","B59","B62","B69","B70","B72","B76","B78","B103","B106",
    "B109","B112","B122","B139","B151","B152","B154","B158"
)

foreach ($addr in $targets) {
    $cell = $ws.Range($addr)
    $val = $cell.Value()
    $words = $val -split ' '
    for ($i = 0; $i -lt $words.Length; $i++) {
        if ($words[$i] -eq 'de') { $words[$i] = 'De' }
        elseif ($words[$i] -eq 'del') { $words[$i] = 'Del' }
        elseif ($words[$i] -eq 'la') { $words[$i] = 'La' }
        elseif ($words[$i] -eq 'los') { $words[$i] = 'Los' }
    }
    $cell.Value = ($words -join ' ')
}

# 3. Delete the footer/metadata rows 170-174 (rows were 170..174, leaving row 169 empty)
$ws.Range("A170:A174").EntireRow.Delete()
